# Estadisticos Matutinos 15 Oct
# Updates the statistics on "Estadisticos 1P", "Estadisticos 2P" and
# "Estadisticos Final" (blanks / reprobados / aprobados / % aprobados /
# promedio for rows 2-8) and fills in the "Rescatables" sheet with the
# list of students who need to retake a subject.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Estadisticos 1P" and "Estadisticos Final" share the same figures.
# ---------------------------------------------------------------------
$blancos    = @(8, 1, 5, 3, 1, 0, 3)
$reprobados = @(0, 9, 7, 4, 10, 6, 6)
$aprobados  = @(20, 18, 21, 24, 10, 29, 12)
$porApro    = @(71.43, 64.29, 63.64, 77.42, 47.62, 82.86, 57.14)
$promedio   = @(7.1, 7.8, 6.8, 7.7, 7.2, 8, 6.9)

foreach ($sheetName in @("Estadisticos 1P", "Estadisticos Final")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt 7; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 4).Value = $blancos[$i]
        $ws.Cells.Item($row, 5).Value = $reprobados[$i]
        $ws.Cells.Item($row, 6).Value = $aprobados[$i]
        $ws.Cells.Item($row, 7).Value = $porApro[$i]
        $ws.Cells.Item($row, 8).Value = $promedio[$i]
    }
}

# ---------------------------------------------------------------------
# 2) "Estadisticos 2P": only the "Reprobados" column changes.
# ---------------------------------------------------------------------
$reprobados2P = @(20, 27, 28, 28, 20, 35, 18)

$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws2P.Cells.Item($row, 5).Value = $reprobados2P[$i]
}

# ---------------------------------------------------------------------
# 3) "Rescatables": list of students with a pending subject to rescue.
# ---------------------------------------------------------------------
$nc       = @(21330051920121, 21330051920121, 20330051920073, 20330051920326, 21330051920118, 21330051920122, 20330051920268)
$paterno  = @("GARCIA", "GARCIA", "CANUTO", "LUNA", "FLORES", "HERNANDEZ", "GONZALEZ")
$materno  = @("RAMIREZ", "RAMIREZ", "MEDINA", "MORALES", "VILLA", "GARCIA", "FLORES")
$nombres  = @("CRISTIAN FERMIN", "CRISTIAN FERMIN", "ISRAEL", "JESUS ANTONIO", "DIEGO", "CLAUDIA", "JESUS HUMBERTO")
$materia  = @("INGLÉS I", "LECTURA, EXPRESIÓN ORAL Y ESCRITA I", "INGLÉS III", "INGLÉS III", "LECTURA, EXPRESIÓN ORAL Y ESCRITA I", "LECTURA, EXPRESIÓN ORAL Y ESCRITA I", "INGLÉS III")
$grupo    = @("1DV", "1DV", "3AEV", "3ASV", "1DV", "1DV", "3APV")
$reprobadas = @(6, 6, 6, 6, 6, 6, 6)

$wsResc = $wb.Worksheets.Item("Rescatables")

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $wsResc.Cells.Item($row, 1).Value = $nc[$i]
}
for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $wsResc.Cells.Item($row, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $wsResc.Cells.Item($row, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $wsResc.Cells.Item($row, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $wsResc.Cells.Item($row, 5).Value = $materia[$i]
}
for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $wsResc.Cells.Item($row, 6).Value = $grupo[$i]
}
for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $wsResc.Cells.Item($row, 7).Value = $reprobadas[$i]
}
